$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 16667802
$ws.Range("I32").Value = 27778378
$ws.Range("J32").Value = 1938.5
$ws.Range("K32").Value = 27778378
$ws.Range("L32").Value = 1938.5
$ws.Range("M32").Value = -27778052
$ws.Range("N32").Value = -2590.5
$ws.Range("H64").Value = 83337496
$ws.Range("I64").Value = 4129.1
$ws.Range("K64").Value = 4129.1
$ws.Range("M64").Value = -3881.1
$ws.Range("H67").Value = 83337496
$ws.Range("I67").Value = 4129.1
$ws.Range("K67").Value = 4129.1
$ws.Range("M67").Value = -3271.1
$ws.Range("H76").Value = 3711.2222
$ws.Range("I76").Value = 3675.125
$ws.Range("K76").Value = 3675.125
$ws.Range("M76").Value = -3360.125
$ws.Range("H79").Value = 3711.2222
$ws.Range("I79").Value = 3675.125
$ws.Range("K79").Value = 3675.125
$ws.Range("M79").Value = -2583.125
$ws.Range("H112").Value = 1697.5217
$ws.Range("J112").Value = 1960.9474
$ws.Range("L112").Value = 5882.8422
$ws.Range("N112").Value = -8098.8422
$ws.Range("H125").Value = 1010.2
$ws.Range("J125").Value = 1014.7143
$ws.Range("L125").Value = 9132.4287
$ws.Range("N125").Value = -14052.4287
$ws.Range("H135").Value = 739.1177
$ws.Range("I135").Value = 469
$ws.Range("J135").Value = 1999.6666
$ws.Range("K135").Value = 4221
$ws.Range("L135").Value = 17996.9994
$ws.Range("M135").Value = -1686
$ws.Range("N135").Value = -23066.9994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 233
$ws.Range("J5").Value = 500
$ws.Range("L5").Value = 500
$ws.Range("N5").Value = -724
$ws.Range("H32").Value = 2290350.8
$ws.Range("I32").Value = 1124786
$ws.Range("K32").Value = 1124786
$ws.Range("M32").Value = -1124499
$ws.Range("H63").Value = 2124.75
$ws.Range("I63").Value = 2099.6667
$ws.Range("K63").Value = 2099.6667
$ws.Range("M63").Value = -1413.6667
$ws.Range("H66").Value = 2124.75
$ws.Range("I66").Value = 2099.6667
$ws.Range("K66").Value = 10498.3335
$ws.Range("M66").Value = -7066.333500000001
$ws.Range("H132").Value = 20000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").Value = 60000
$ws.Range("N132").Value = -65060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 233
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -730
$ws.Range("H99").Value = 58029.5
$ws.Range("I99").Value = 73759.42999999999
$ws.Range("J99").Value = 2974.75
$ws.Range("K99").Value = 73759.42999999999
$ws.Range("L99").Value = 2974.75
$ws.Range("M99").Value = -72261.42999999999
$ws.Range("N99").Value = -5970.75
$ws.Range("H103").Value = 37006.875
$ws.Range("J103").Value = 37006.875
$ws.Range("L103").Value = 37006.875
$ws.Range("N103").Value = -39350.875
$ws.Range("H107").Value = 2263834.8
$ws.Range("I107").Value = 2653852.5
$ws.Range("K107").Value = 2653852.5
$ws.Range("M107").Value = -2651932.5
$ws.Range("H134").Value = 2413.5715
$ws.Range("I134").Value = 1979
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 5937
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -3402
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H134").Value = 4053.1904
$ws.Range("I134").Value = 4053.1904
$ws.Range("K134").Value = 12159.5712
$ws.Range("M134").Value = -9624.5712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4214.9585
$ws.Range("J134").Value = 4980.1577
$ws.Range("L134").Value = 14940.4731
$ws.Range("N134").Value = -25080.4731

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 42.25
$ws.Range("I2").Value = 31.333334
$ws.Range("J2").Value = 75
$ws.Range("K2").Value = 31.333334
$ws.Range("L2").Value = 75
$ws.Range("M2").Value = 81.66666599999999
$ws.Range("N2").Value = -301
$ws.Range("H112").Value = 98999.5
$ws.Range("J112").Value = 98999.5
$ws.Range("L112").Value = 98999.5
$ws.Range("N112").Value = -101215.5
$ws.Range("H122").Value = 4052296.5
$ws.Range("I122").Value = 5919508.5
$ws.Range("J122").Value = 6670.5
$ws.Range("K122").Value = 17758525.5
$ws.Range("L122").Value = 20011.5
$ws.Range("M122").Value = -17756075.5
$ws.Range("N122").Value = -24911.5
$ws.Range("H132").Value = 62503216
$ws.Range("I132").Value = 200002990
$ws.Range("J132").Value = 3317.2727
$ws.Range("K132").Value = 600008970
$ws.Range("L132").Value = 9951.8181
$ws.Range("M132").Value = -600006440
$ws.Range("N132").Value = -15011.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1475
$ws.Range("I22").Value = 1475
$ws.Range("K22").Value = 1475
$ws.Range("M22").Value = -1180
$ws.Range("H27").Value = 1475
$ws.Range("I27").Value = 1475
$ws.Range("K27").Value = 1475
$ws.Range("M27").Value = -1368
$ws.Range("H40").Value = 41000.2
$ws.Range("I40").Value = 41000.2
$ws.Range("K40").Value = 41000.2
$ws.Range("M40").Value = -40864.2
$ws.Range("H46").Value = 3185.7144
$ws.Range("I46").Value = 2466.6667
$ws.Range("K46").Value = 2466.6667
$ws.Range("M46").Value = -2278.6667
$ws.Range("H100").Value = 3031.0833
$ws.Range("I100").Value = 1782
$ws.Range("K100").Value = 1782
$ws.Range("M100").Value = -1241
$ws.Range("H122").Value = 5458.2563
$ws.Range("I122").Value = 3769.1667
$ws.Range("K122").Value = 11307.5001
$ws.Range("M122").Value = -8857.500100000001
$ws.Range("H132").Value = 4709.1924
$ws.Range("I132").Value = 4673.1177
$ws.Range("K132").Value = 14019.3531
$ws.Range("M132").Value = -11489.3531
$ws.Range("H133").Value = 83444
$ws.Range("J133").Value = 83444
$ws.Range("L133").Value = 83444
$ws.Range("N133").Value = -88504

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 36928.6
$ws.Range("J101").Value = 36928.6
$ws.Range("L101").Value = 36928.6
$ws.Range("N101").Value = -43418.6
$ws.Range("H122").Value = 15626999
$ws.Range("I122").Value = 2132.1333
$ws.Range("K122").Value = 6396.3999
$ws.Range("M122").Value = -3946.3999
$ws.Range("H126").Value = 2015.6451
$ws.Range("I126").Value = 1820.2
$ws.Range("K126").Value = 5460.6
$ws.Range("M126").Value = -2990.6
$ws.Range("H132").Value = 2168.348
$ws.Range("I132").Value = 1893.2632
$ws.Range("K132").Value = 5679.7896
$ws.Range("M132").Value = -3149.7896
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
